$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting existing rows 29..107 down to 30..108
$ws.Rows("29:29").Insert()

# Populate the new row 29 with the new weekly price record
$ws.Range("A29").Value = 11
$ws.Range("B29").Value = "Vega Monumental Concepción"
$ws.Range("C29").Value = "Bíobío"
$ws.Range("D29").Value = 44624
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = 100112043
$ws.Range("G29").Value = "Pepino ensalada"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 270
$ws.Range("K29").Value = 14000
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = 14444
$ws.Range("N29").Value = "$/caja 60 unidades"
$ws.Range("O29").Value = "Región Metropolitana"
$ws.Range("P29").Value = 241
$ws.Range("Q29").Value = 60
$ws.Range("R29").Value = "Hortaliza"
